$d = $word.ActiveDocument

# 1. Insert "CS " before "capstone sequence." in the "These courses provide students..." paragraph.
$range1 = $d.Content
$range1.Find.Execute("complementing the capstone sequence.", $true, $false, $false, $false, $false, $true, 1, $false, "complementing the CS capstone sequence.", 2)

# 2. Replace "enables" with "provides" in "The capstone sequence..." paragraph.
$range2 = $d.Content
$range2.Find.Execute("enables", $true, $false, $false, $false, $false, $true, 1, $false, "provides", 2)
